$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.260.54"
$ws.Range("D3").Value = "'1.907.15"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'307.84"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.5252"
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("D8").Value = "'0.3791"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("D9").Value = "'0.07268"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'21.32"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("D11").Value = "'0.9007"
$ws.Range("D12").Value = "'0.08139"
$ws.Range("E12").Value = "  +8.75%  "
$ws.Range("D13").Value = "'1.911.22"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "'95.43"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "'5.295"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'0.000008632"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "'14.50"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'27.323.55"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'5.067"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "'2.149.96"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").Value = "'6.461"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'2.312"
$ws.Range("E25").Value = "  +11.10%  "
$ws.Range("D26").Value = "'146.30"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'1.747"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'18.20"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "'114.92"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "'4.996"
$ws.Range("E30").Value = "  +6.29%  "
$ws.Range("D31").Value = "'4.815"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").Value = "'0.09229"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'0.8064"
$ws.Range("E33").Value = "  +7.93%  "
$ws.Range("D34").Value = "'0.05063"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +8.22%  "
$ws.Range("D36").Value = "'2.997"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "'3.327"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").Value = "'2.582"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("D39").Value = "'0.5744"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "'0.01988"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'1.079"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "'119.66"
$ws.Range("D43").Value = "'6.634"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").Value = "'8.975"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").Value = "'0.4860"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "'10.26"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("D50").Value = "'37.65"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("D51").Value = "'63.81"
$ws.Range("E51").Value = "  +1.41%  "
